# Actualizacion automatica del tracker
# - Completa resultado (G) / profit (H) de partidos ya listados (filas 43, 45, 46)
# - Agrega los nuevos partidos del 2025-08-05 al final de la tabla (filas 55-58)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resultados de partidos ya registrados ---
$ws.Range("G43").Value = "Fallo"
$ws.Range("H43").Value = -1

$ws.Range("G45").Value = "Acierto"
$ws.Range("H45").Value = 1.63

$ws.Range("G46").Value = "Fallo"
$ws.Range("H46").Value = -1

# --- Nuevos partidos trackeados ---
$newRows = @(
    @(14349613, "2025-08-05", "Francesco Maestrelli", "Maximus Jones", "Gana Maximus Jones", 3.5),
    @(14349730, "2025-08-05", "Zdenek Kolar", "Zsombor Piros", "Gana Zdenek Kolar", 3.4),
    @(14349723, "2025-08-05", "Gonzalo Bueno", "Ryan Nijboer", "Gana Ryan Nijboer", 2.75),
    @(14349602, "2025-08-05", "Mats Rosenkranz", "Clement Tabur", "Gana Mats Rosenkranz", 3.25)
)

$row = 55
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]

    # Forzar texto para que la fecha ISO no se auto-convierta a numero de serie
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $data[1]

    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]

    # resultado / profit aun no disponibles para estos partidos
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = ""

    $row = $row + 1
}
